# Update the RandomForest-imputed result values for the terrestrial_mammals
# BCE/20/seed3 combination ("Update Name of Algo" re-run refreshed numbers).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value = 16.22890000000001
$ws.Range("E4").Value = 16.36189999999998
$ws.Range("C7").Value = -12.76009999999999
$ws.Range("B8").Value = 5.710199999999997
$ws.Range("B10").Value = 6.603799999999995
$ws.Range("E11").Value = 15.88590000000001
$ws.Range("B12").Value = 4.542599999999998
$ws.Range("C14").Value = -13.8878
$ws.Range("E14").Value = 16.36770000000001
$ws.Range("C15").Value = -14.03789999999999
$ws.Range("B18").Value = 7.292600000000005
$ws.Range("C18").Value = -11.34879999999999
$ws.Range("E18").Value = 18.05980000000002
$ws.Range("E19").Value = 16.55130000000001
$ws.Range("C20").Value = -12.17
$ws.Range("E21").Value = 16.7041
$ws.Range("B25").Value = 5.452400000000001
$ws.Range("E27").Value = 16.7118
$ws.Range("C29").Value = -11.3861
$ws.Range("C30").Value = -13.3429
$ws.Range("C31").Value = -13.4135
$ws.Range("E31").Value = 16.44410000000001
$ws.Range("C35").Value = -11.4191
$ws.Range("B37").Value = 8.966800000000001
$ws.Range("E38").Value = 16.42849999999999
$ws.Range("C40").Value = -12.7393
$ws.Range("E42").Value = 16.2223
$ws.Range("C44").Value = -12.4576
$ws.Range("E44").Value = 16.82749999999999
$ws.Range("E47").Value = 16.7214
$ws.Range("C50").Value = -13.30769999999999
$ws.Range("C54").Value = -13.2054
$ws.Range("B55").Value = 6.180499999999995
$ws.Range("E56").Value = 16.5564
$ws.Range("E58").Value = 16.03430000000002
$ws.Range("E65").Value = 17.3435
$ws.Range("B68").Value = 6.462499999999998
$ws.Range("C68").Value = -11.7901
$ws.Range("E73").Value = 17.4129
$ws.Range("C76").Value = -12.3811
$ws.Range("B77").Value = 9.745400000000004
$ws.Range("B78").Value = 9.519000000000004
$ws.Range("B79").Value = 9.351700000000003
$ws.Range("B80").Value = 9.449299999999999
$ws.Range("B81").Value = 6.496500000000001
$ws.Range("B82").Value = 4.959300000000002
$ws.Range("B84").Value = 6.9815
$ws.Range("C87").Value = -13.74659999999999
$ws.Range("C88").Value = -12.6225
$ws.Range("E90").Value = 16.52069999999999
$ws.Range("C92").Value = -10.8404
$ws.Range("E92").Value = 18.44140000000002
$ws.Range("E94").Value = 19.09850000000002
$ws.Range("E95").Value = 18.31560000000002
$ws.Range("C96").Value = -13.13350000000001
$ws.Range("C98").Value = -11.94089999999999
$ws.Range("B101").Value = 9.091999999999997
$ws.Range("C101").Value = -12.35180000000001
$ws.Range("E101").Value = 16.43070000000001
$ws.Range("B102").Value = 8.349100000000009
$ws.Range("C102").Value = -13.32240000000001
